# Update row 9 (year 2025) figures in the faturamento_anual worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3783252.78
$ws.Range("C9").Value = 593895.24
$ws.Range("D9").Value = 4377148.02
$ws.Range("E9").Value = 13.56808673790291
$ws.Range("F9").Value = 86.4319132620971
$ws.Range("G9").Value = -42.6029863521874
$ws.Range("H9").Value = -31.67968415506432
$ws.Range("I9").Value = 38197
$ws.Range("J9").Value = 1623
$ws.Range("K9").Value = 39820
$ws.Range("L9").Value = 27507
$ws.Range("M9").Value = 159.1285134692987
$ws.Range("N9").Value = 8.640277231993942
